# Update Legs on Tester Higher Speed Sesi 1!
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("D4").Value = 1580
$ws.Range("G4").Value = 2050

# Row 5
$ws.Range("C5").Value = 1520
$ws.Range("D5").Value = 1500
$ws.Range("F5").Value = 2080
$ws.Range("G5").Value = 2030

# Row 6
$ws.Range("C6").Value = 1350
$ws.Range("G6").Value = 2150

# Row 7
$ws.Range("C7").Value = 1410
$ws.Range("D7").Value = 1380

# Row 8
$ws.Range("C8").Value = 1600
$ws.Range("D8").Value = 1250

# Update the active cell selection to match the saved view state
$ws.Range("G15").Select()
